# Added save generated report test case
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: the existing "Generate report" test case gets a proper Id ---
$ws.Range("B2").Value = "T_C_105.1"

# --- Row 3: "Save generated report" test case gets its sub-level Id ---
$ws.Range("B3").Value = "T_C_105.2"

# --- Row 4 (new): sub-level test case for the "Save generated report" task ---
# Clone the border-only format used by A3:D3 across A4:G4, and the bold
# green "Pass" status format used by H3 onto H4, so the new row matches the
# rest of the table exactly.
$ws.Range("A3").Copy()
$ws.Range("A4:G4").PasteSpecial(-4122)
$ws.Range("H3").Copy()
$ws.Range("H4").PasteSpecial(-4122)

$ws.Range("D4").Value = "Report button is clicked before any cultivation "
$ws.Range("E4").Value = "Medium"
$ws.Range("B4").Value = "T_C_105.3"
$ws.Range("F4").Value = "Tester A"
$ws.Range("G4").NumberFormat = "mm-dd-yy"
$ws.Range("G4").Value = "1/12/2017"
$ws.Range("H4").Value = "Pass"

# --- Row 5: blank row, but still boxed with the same thin border ---
$ws.Range("A4").Copy()
$ws.Range("A5:H5").PasteSpecial(-4122)
$ws.Range("A5:H5").ClearContents()

$excel.CutCopyMode = 0

# --- Selection moves to D7 ---
$ws.Range("D7").Select()
